$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "Sheet1"

# Update / add the login data rows
$data = @(
    @("admin", "admin"),
    @("andre@pais", "Andre1234"),
    @("andre@silva", "Andre1234"),
    @("eduardo@pais", "Eduardo1234"),
    @("renaro@gmail.com", "Renato1234"),
    @("andre@1234", "Andre-123"),
    @("andre@1234", "andre 123"),
    @("andre@pais", "andre1234"),
    @("andre@pais", "Andre1234"),
    @("andre@pais", "Andre123"),
    @("pao@manteiga", "Pao1234567"),
    @("manteiga@pao", "Pao-1234567"),
    @("man@man", "Man1234567"),
    @("andre@pais", "Andre1234")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
